$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended right after the current last row (69) of the ticket log.
$newRow = 70
$lastRow = 69

$targetRange = $ws.Range("A" + $newRow + ":I" + $newRow)

# Force the new cells to be treated as plain text first, so values that look
# like dates/times (e.g. "2024-05-20", "10:51:43") are stored as literal
# strings instead of being auto-converted to date/time serial numbers.
$targetRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2024-05-20"
$ws.Cells.Item($newRow, 2).Value = "10:51:43"
$ws.Cells.Item($newRow, 3).Value = "-"
$ws.Cells.Item($newRow, 4).Value = "-"
$ws.Cells.Item($newRow, 5).Value = "Detección de sealling mal puesto"
$ws.Cells.Item($newRow, 6).Value = "-"
$ws.Cells.Item($newRow, 7).Value = "-"
$ws.Cells.Item($newRow, 8).Value = "10:52:50"
$ws.Cells.Item($newRow, 9).Value = "0:01:07"

# Re-apply the formatting of the previous row onto the new one so the added
# row keeps the same (default) cell style as the rest of the sheet, instead
# of the explicit "Text" number format applied above.
$ws.Range("A" + $lastRow + ":I" + $lastRow).Copy()
$targetRange.PasteSpecial(-4122)
